$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("ass sizes are smaller, which would definitely affect student performance", $true, $false, $false, $false, $false, $true, 1, $false, "ass sizes are smaller, which would definitely affect student performance.", 2)
